# Apply the 2026-01-05 timesheet edits: simulator full-month coverage,
# persist logs, fix employees.
$wb = $excel.ActiveWorkbook

$wsTime = $wb.Worksheets.Item("Weekly Timesheet")
$wsSchema = $wb.Worksheets.Item("Jason Schema")

# New employee id for Boban Abbate
$newEmpId = "emp_ga4rqytu"

# ---------------------------------------------------------------------
# Weekly Timesheet sheet: columns A Date | B Client | C Hours | D Type |
#                          E Rate | F Total
# ---------------------------------------------------------------------

# Row 2 - 2026-01-05 becomes a PTO day for PTO (client/type now "PTO")
$wsTime.Range("B2").Value = "PTO"
$wsTime.Range("C2").Value = 6.5
$wsTime.Range("D2").Value = "PTO"
$wsTime.Range("E2").Value = 110
$wsTime.Range("F2").Value = 715

# Row 3 - 2026-01-06 worked by Muncey
$wsTime.Range("B3").Value = "Muncey"
$wsTime.Range("C3").Value = 7
$wsTime.Range("D3").Value = "Regular"
$wsTime.Range("E3").Value = 110
$wsTime.Range("F3").Value = 770

# Row 4 - 2026-01-07 worked by Moulton
$wsTime.Range("B4").Value = "Moulton"
$wsTime.Range("C4").Value = 6
$wsTime.Range("D4").Value = "Regular"
$wsTime.Range("E4").Value = 110
$wsTime.Range("F4").Value = 660

# Row 5 - 2026-01-08 worked by Regan
$wsTime.Range("B5").Value = "Regan"
$wsTime.Range("C5").Value = 6
$wsTime.Range("D5").Value = "Regular"
$wsTime.Range("E5").Value = 110
$wsTime.Range("F5").Value = 660

# Row 6 - 2026-01-09 worked by Hendricks
$wsTime.Range("B6").Value = "Hendricks"
$wsTime.Range("C6").Value = 6.5
$wsTime.Range("D6").Value = "Regular"
$wsTime.Range("E6").Value = 110
$wsTime.Range("F6").Value = 715

# Row 8 - SUBTOTAL: hours now total 32, and the total pay is 3520
$wsTime.Range("C8").Value = 32
$wsTime.Range("D8").Value = "Reg: 32 / OT: 0"
$wsTime.Range("F8").Value = 3520

# Row 11 - HOURLY SUBTOTAL total
$wsTime.Range("F11").Value = 3520

# Row 13 - GRAND TOTAL total
$wsTime.Range("F13").Value = 3520

# ---------------------------------------------------------------------
# Jason Schema sheet: columns A Employee | B Employee ID | C Date |
#                      D Client | E Hours | F Rate | G Total | H Type |
#                      I Notes
# ---------------------------------------------------------------------

# Row 2 - 2026-01-05 PTO
$wsSchema.Range("B2").Value = $newEmpId
$wsSchema.Range("D2").Value = "PTO"
$wsSchema.Range("E2").Value = 6.5
$wsSchema.Range("F2").Value = 110
$wsSchema.Range("G2").Value = 715
$wsSchema.Range("H2").Value = "PTO"
$wsSchema.Range("I2").Value = "PTO"

# Row 3 - 2026-01-06 Muncey
$wsSchema.Range("B3").Value = $newEmpId
$wsSchema.Range("D3").Value = "Muncey"
$wsSchema.Range("E3").Value = 7
$wsSchema.Range("F3").Value = 110
$wsSchema.Range("G3").Value = 770
$wsSchema.Range("H3").Value = "Regular"

# Row 4 - 2026-01-07 Moulton
$wsSchema.Range("B4").Value = $newEmpId
$wsSchema.Range("D4").Value = "Moulton"
$wsSchema.Range("E4").Value = 6
$wsSchema.Range("F4").Value = 110
$wsSchema.Range("G4").Value = 660
$wsSchema.Range("H4").Value = "Regular"

# Row 5 - 2026-01-08 Regan
$wsSchema.Range("B5").Value = $newEmpId
$wsSchema.Range("D5").Value = "Regan"
$wsSchema.Range("E5").Value = 6
$wsSchema.Range("F5").Value = 110
$wsSchema.Range("G5").Value = 660
$wsSchema.Range("H5").Value = "Regular"

# Row 6 - 2026-01-09 Hendricks
$wsSchema.Range("B6").Value = $newEmpId
$wsSchema.Range("D6").Value = "Hendricks"
$wsSchema.Range("E6").Value = 6.5
$wsSchema.Range("F6").Value = 110
$wsSchema.Range("G6").Value = 715
$wsSchema.Range("H6").Value = "Regular"
